# Experiment Report.xlsx update — "Updating code for image prediction"
#
# 1) I18's "best so far" note loses the parenthetical and gains a trailing space.
# 2) A new experiment row (34 / output-34) is appended as row 36.
# 3) I1 header note gets a more detailed explanation.
# 4) Selection moves to I18 (was I35).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. Update the row 18 commentary text (I18) ---
$ws.Range("I18").Value = "Valid output `nmacro corelation(54-0%) micro corelation(94-63%)"

# --- 2. Append new experiment row 36 ---
$ws.Range("A36").Value = 34
$ws.Range("B36").Value = $false
$ws.Range("C36").Value = 30
$ws.Range("D36").Value = 0.2
$ws.Range("E36").Value = -1
$ws.Range("F36").Value = "32x32"
$ws.Range("F36").WrapText = $true
$ws.Range("G36").Value = "64x64"
$ws.Range("H36").Value = "output-34"
$ws.Range("I36").Value = "valid output, circle-rectangle macro correlation is high`nmacro corelation(57-5%) micro corelation(94-61%)"
$ws.Range("I36").WrapText = $true
$ws.Rows.Item(36).RowHeight = 36.75

# --- 3. Update the header result-column description (I1) ---
$ws.Range("I1").Value = "Result(micro and macro corelation where the percentages refer to the maximun and minimum values from all the shapes)"

# --- 4. Move the active selection to I18 ---
$ws.Range("I18").Select()
